# Add a new "bread" worksheet at the end of the workbook (after "noodles"),
# populate it with property/value rows (mirroring the existing food-group
# sheets), and leave it as the active/selected sheet - matching the
# "added more csv files" commit.

$wb = $excel.ActiveWorkbook

# The previously-active sheet ("noodles") ends up with its whole data range
# selected (and loses the single-cell selection it had below its data).
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$lastSheet.Range("A1:B6").Select()

# Insert the new sheet right after the current last worksheet ("noodles")
# so it lands at the end of the tab strip.
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "bread"

# property / value header
$ws.Range("A1").Value = "property"
$ws.Range("B1").Value = "value"

# name
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "noodles"

# healthy (stored as literal text "true", not a boolean - avoid Excel's
# auto-conversion of the bare word "true" into a Boolean TRUE cell)
$ws.Range("A3").Value = "healthy"
$ws.Range("B3").Formula = "=""true"""
$ws.Range("B3").Value = $ws.Range("B3").Value

# food super group
$ws.Range("A4").Value = "food super group"
$ws.Range("B4").Value = "carbs"

# gluten-free substitute
$ws.Range("A5").Value = "gluten-free substitute"
$ws.Range("B5").Value = "pass"

# style
$ws.Range("A6").Value = "style"
$ws.Range("B6").Value = "pass"

# Match the author's final selection on the new sheet (cell below the data).
$ws.Range("B7").Select()
